# Weekly update: a new Mango price record for Terminal Hortofrutícola Agro
# Chillán is inserted as the new row 48, pushing the existing rows 48:100
# down to 49:101 (dimension grows from A1:T100 to A1:T101).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 48; everything below
# (old rows 48:100) shifts down to 49:101, carrying its formatting along.
$ws.Rows(48).Insert()

# Populate the newly inserted row 48 with this week's record.
$ws.Range("A48").Value = 7
$ws.Range("B48").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C48").Value = "Ñuble"
$ws.Range("D48").Value = 44994
$ws.Range("E48").Value = 16
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100108
$ws.Range("H48").Value = "Tropicales y subtropicales"
$ws.Range("I48").Value = 100108002
$ws.Range("J48").Value = "Mango"
$ws.Range("K48").Value = "Sin especificar"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 30
$ws.Range("N48").Value = 8000
$ws.Range("O48").Value = 8000
$ws.Range("P48").Value = 8000
$ws.Range("Q48").Value = "`$/bandeja 4 kilos"
$ws.Range("R48").Value = "Perú"
$ws.Range("S48").Value = 2000
$ws.Range("T48").Value = 4
